# Natmi following Dr Hou advice
# Recompute the ligand/receptor edge table (Tnfsf13b -> Tnfrsf13b) adding the
# EC and sC sending clusters alongside the existing FAPs/M2 clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> M2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf13b"
$ws.Range("C2").Value = "Tnfrsf13b"
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4101656666666667
$ws.Range("H2").Value = 1.230497
$ws.Range("I2").Value = 0.07383296380759893
$ws.Range("J2").Value = 0.07383296380759893
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.557669333333333
$ws.Range("N2").Value = 28.673008
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 3.920227813886223
$ws.Range("R2").Value = 35.28205032497601
$ws.Range("S2").Value = 0.07383296380759893
$ws.Range("T2").Value = 0.07383296380759893

# Row 3: FAPs -> M2
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf13b"
$ws.Range("C3").Value = "Tnfrsf13b"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.119982333333334
$ws.Range("H3").Value = 12.359947
$ws.Range("I3").Value = 0.7416283985372097
$ws.Range("J3").Value = 0.7416283985372097
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.557669333333333
$ws.Range("N3").Value = 28.673008
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 39.37742880117511
$ws.Range("R3").Value = 354.3968592105761
$ws.Range("S3").Value = 0.7416283985372097
$ws.Range("T3").Value = 0.7416283985372097

# Row 4 (new): M2 -> M2
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Tnfsf13b"
$ws.Range("C4").Value = "Tnfrsf13b"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7080713333333333
$ws.Range("H4").Value = 2.124214
$ws.Range("I4").Value = 0.1274582671730162
$ws.Range("J4").Value = 0.1274582671730162
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.557669333333333
$ws.Range("N4").Value = 28.673008
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 6.767511668412443
$ws.Range("R4").Value = 60.90760501571199
$ws.Range("S4").Value = 0.1274582671730162
$ws.Range("T4").Value = 0.1274582671730162

# Row 5 (new): sCs -> M2
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Tnfsf13b"
$ws.Range("C5").Value = "Tnfrsf13b"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3170996666666667
$ws.Range("H5").Value = 0.951299
$ws.Range("I5").Value = 0.05708037048217513
$ws.Range("J5").Value = 0.05708037048217512
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.557669333333333
$ws.Range("N5").Value = 28.673008
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 3.030733759710222
$ws.Range("R5").Value = 27.276603837392
$ws.Range("S5").Value = 0.05708037048217513
$ws.Range("T5").Value = 0.05708037048217512
